$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A48").Value = "III-2021"
$ws.Range("B48").Value = 4182.8
$ws.Range("C48").Value = 3797.3
$ws.Range("D48").Value = 385.5
